# TimeScheduleConfig.xlsx update:
# - EBOVersion on the Configuration sheet bumps from 6.0.4.90 to 5.0.3.117
# - The Configuration sheet becomes the active tab/selection (was on Entries)

$wb = $excel.ActiveWorkbook

$entries = $wb.Worksheets.Item("Entries")
$config  = $wb.Worksheets.Item("Configuration")

# Bump the EBOVersion value in Configuration!B6
$config.Range("B6").Value = "5.0.3.117"

# Move the active sheet/selection from Entries to Configuration (B7)
$config.Activate()
$config.Range("B7").Select()
